$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 283.48276
$ws.Range("I28").Value = 269.54166
$ws.Range("J28").Value = 350.4
$ws.Range("K28").Value = 269.54166
$ws.Range("L28").Value = 350.4
$ws.Range("M28").Value = 215.45834
$ws.Range("N28").Value = -1320.4
$ws.Range("H39").Value = 103.833336
$ws.Range("I39").Value = 51.266666
$ws.Range("J39").Value = 366.66666
$ws.Range("K39").Value = 153.799998
$ws.Range("L39").Value = 1099.99998
$ws.Range("M39").Value = 142.200002
$ws.Range("N39").Value = -1691.99998
$ws.Range("H51").Value = 2960
$ws.Range("I51").Value = 1900
$ws.Range("K51").Value = 1900
$ws.Range("M51").Value = -1416
$ws.Range("H70").Value = 51314.5
$ws.Range("I70").Value = 251075
$ws.Range("J70").Value = 1374.375
$ws.Range("K70").Value = 753225
$ws.Range("L70").Value = 4123.125
$ws.Range("M70").Value = -752955
$ws.Range("N70").Value = -4663.125
$ws.Range("H73").Value = 51314.5
$ws.Range("I73").Value = 251075
$ws.Range("J73").Value = 1374.375
$ws.Range("K73").Value = 753225
$ws.Range("L73").Value = 4123.125
$ws.Range("M73").Value = -752289
$ws.Range("N73").Value = -5995.125
$ws.Range("H103").Value = 781.2857
$ws.Range("I103").Value = 773.8
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 2321.4
$ws.Range("L103").Value = 2400
$ws.Range("M103").Value = -1735.4
$ws.Range("N103").Value = -3572
$ws.Range("H129").Value = 1038.1569
$ws.Range("J129").Value = 1045.8368
$ws.Range("L129").Value = 3137.5104
$ws.Range("N129").Value = -13137.5104
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 66487.5
$ws.Range("J136").Value = 66487.5
$ws.Range("L136").Value = 66487.5
$ws.Range("N136").Value = -76687.5

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 992
$ws.Range("I2").Value = 779
$ws.Range("J2").Value = 1577.75
$ws.Range("K2").Value = 779
$ws.Range("L2").Value = 1577.75
$ws.Range("M2").Value = -666
$ws.Range("N2").Value = -1803.75
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H32").Value = 10165.396
$ws.Range("I32").Value = 10654.576
$ws.Range("J32").Value = 2950
$ws.Range("K32").Value = 10654.576
$ws.Range("L32").Value = 2950
$ws.Range("M32").Value = -10367.576
$ws.Range("N32").Value = -3524
$ws.Range("H116").Value = 992
$ws.Range("I116").Value = 779
$ws.Range("J116").Value = 1577.75
$ws.Range("K116").Value = 779
$ws.Range("L116").Value = 1577.75
$ws.Range("M116").Value = 1515
$ws.Range("N116").Value = -6165.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 992
$ws.Range("I3").Value = 779
$ws.Range("J3").Value = 1577.75
$ws.Range("K3").Value = 779
$ws.Range("L3").Value = 1577.75
$ws.Range("M3").Value = -665
$ws.Range("N3").Value = -1805.75
$ws.Range("H32").Value = 50000
$ws.Range("I32").Value = 50000
$ws.Range("K32").Value = 50000
$ws.Range("M32").Value = -49616
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H107").Value = 28187.422
$ws.Range("I107").Value = 35336.934
$ws.Range("J107").Value = 1376.75
$ws.Range("K107").Value = 35336.934
$ws.Range("L107").Value = 1376.75
$ws.Range("M107").Value = -33416.934
$ws.Range("N107").Value = -5216.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1791.9166
$ws.Range("I16").Value = 1671.5714
$ws.Range("J16").Value = 1960.4
$ws.Range("K16").Value = 1671.5714
$ws.Range("L16").Value = 1960.4
$ws.Range("M16").Value = -1384.5714
$ws.Range("N16").Value = -2534.4
$ws.Range("H113").Value = 1791.9166
$ws.Range("I113").Value = 1671.5714
$ws.Range("J113").Value = 1960.4
$ws.Range("K113").Value = 1671.5714
$ws.Range("L113").Value = 1960.4
$ws.Range("M113").Value = 498.4286
$ws.Range("N113").Value = -6300.4
$ws.Range("H122").Value = 2793.862
$ws.Range("I122").Value = 2733.7778
$ws.Range("J122").Value = 2892.182
$ws.Range("K122").Value = 8201.3334
$ws.Range("L122").Value = 8676.545999999998
$ws.Range("M122").Value = -5751.3334
$ws.Range("N122").Value = -13576.546

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 7000
$ws.Range("J35").Value = 7000
$ws.Range("L35").Value = 21000
$ws.Range("N35").Value = -21576
$ws.Range("H120").Value = 8151.5
$ws.Range("I120").Value = 6324.222
$ws.Range("K120").Value = 18972.666
$ws.Range("M120").Value = -14134.666
$ws.Range("H131").Value = 14287134
$ws.Range("I131").Value = 480
$ws.Range("J131").Value = 15386107
$ws.Range("K131").Value = 1440
$ws.Range("L131").Value = 46158321
$ws.Range("M131").Value = 3600
$ws.Range("N131").Value = -46168401
$ws.Range("H134").Value = 4430.0225
$ws.Range("I134").Value = 1709.375
$ws.Range("K134").Value = 5128.125
$ws.Range("M134").Value = -58.125

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1622.4667
$ws.Range("J122").Value = 3750
$ws.Range("L122").Value = 11250
$ws.Range("N122").Value = -16150
$ws.Range("H125").Value = 49065
$ws.Range("J125").Value = 49065
$ws.Range("L125").Value = 49065
$ws.Range("N125").Value = -53985

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1066.3636
$ws.Range("I22").Value = 550
$ws.Range("J22").Value = 1970
$ws.Range("K22").Value = 550
$ws.Range("L22").Value = 1970
$ws.Range("M22").Value = -255
$ws.Range("N22").Value = -2560
$ws.Range("H27").Value = 1066.3636
$ws.Range("I27").Value = 550
$ws.Range("J27").Value = 1970
$ws.Range("K27").Value = 550
$ws.Range("L27").Value = 1970
$ws.Range("M27").Value = -443
$ws.Range("N27").Value = -2184
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 3972.08
$ws.Range("I100").Value = 4268.5264
$ws.Range("J100").Value = 3033.3333
$ws.Range("K100").Value = 4268.5264
$ws.Range("L100").Value = 3033.3333
$ws.Range("M100").Value = -3727.5264
$ws.Range("N100").Value = -4115.3333
$ws.Range("H122").Value = 13893800
$ws.Range("I122").Value = 19234954
$ws.Range("K122").Value = 57704862
$ws.Range("M122").Value = -57702412

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 13753.25
$ws.Range("I14").Value = 15004
$ws.Range("J14").Value = 10001
$ws.Range("K14").Value = 15004
$ws.Range("L14").Value = 10001
$ws.Range("M14").Value = -14836
$ws.Range("N14").Value = -10337
$ws.Range("H96").Value = 1900
$ws.Range("J96").Value = 1900
$ws.Range("L96").Value = 1900
$ws.Range("N96").Value = -4646
$ws.Range("H122").Value = 7815002.5
$ws.Range("I122").Value = 11365594
$ws.Range("J122").Value = 3702
$ws.Range("K122").Value = 34096782
$ws.Range("L122").Value = 11106
$ws.Range("M122").Value = -34094332
$ws.Range("N122").Value = -16006
$ws.Range("H126").Value = 17294.273
$ws.Range("I126").Value = 17294.273
$ws.Range("K126").Value = 51882.819
$ws.Range("M126").Value = -49412.819
$ws.Range("H136").Value = 1638.762
$ws.Range("I136").Value = 1807.25
$ws.Range("J136").Value = 1099.6
$ws.Range("K136").Value = 5421.75
$ws.Range("L136").Value = 3298.8
$ws.Range("M136").Value = -2871.75
$ws.Range("N136").Value = -8398.799999999999
